# Generate Report for Handoff
# Update the status/handoff info for the "6c476987-2e90-41aa-b3d3-52976c12f200.md" file
# row, across the Overview, zh-cn and de-de sheets, to reflect it is now ready for handoff.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-03-09 18:41:38"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-03-09 18:41:42"
